# QPEDataClean.xlsx — "added boxplot analysis of QPE"
#
# Content-level changes captured by the diff:
#   1. Header cell C1 is renamed from "Q# loop in Q#" to
#      "Q# loop inside of Q#" (Excel drops the now-unused shared string and
#      appends the edited text at the end of the shared-string table, which
#      is why the <v> indices for C1/D1/E1 shuffle in the XML — the visible
#      text of D1/E1 is unchanged).
#   2. The user's selection moves from A1:E11 (active cell E11) to the
#      single cell I13.
#
# (The remaining hunks in the diff — absPath username, revisionPtr ids,
# window geometry, the default-row-height/dyDescent metrics and the
# sub-pixel "bestFit" column-width jitter — are host-environment artifacts
# written by the authoring machine's Excel build, not deliberate
# spreadsheet edits, and are not exposed through the Excel object model.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Q# loop inside of Q#"

$ws.Range("I13").Select() | Out-Null
